# Fix typo in framework: "Balance sheet total" -> "Balance Sheet Total"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7, column D holds the "Balance sheet total" field-name cell.
$ws.Range("D7").Value = "Balance Sheet Total"

# Mirror the author's final cursor position after the edit.
$null = $ws.Range("D8").Select()
